$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D19").Value = "隔壁学生"
$ws.Range("E19").Value = "Y"
$ws.Range("D20").Value = "隔壁学生"
$ws.Range("E20").Value = "Y"
$ws.Range("D21").Value = "隔壁学生"
$ws.Range("E21").Value = "Y"
$ws.Range("D22").Value = "隔壁学生"
$ws.Range("E22").Value = "Y"
$ws.Range("D23").Value = "隔壁学生"
$ws.Range("E23").Value = "Y"
$ws.Range("D24").Value = "隔壁学生"
$ws.Range("E24").Value = "Y"
$ws.Range("D25").Value = "隔壁学生"
$ws.Range("E25").Value = "Y"
$ws.Range("D26").Value = "隔壁学生"
$ws.Range("E26").Value = "Y"
$ws.Range("D27").Value = "隔壁学生"
$ws.Range("E27").Value = "Y"
$ws.Range("E28").Value = "Y"
$ws.Range("E30").Value = "Y"
$ws.Range("E31").Value = "Y"
$ws.Range("E32").Value = "Y"
$ws.Range("E33").Value = "Y"
$ws.Range("D39").Value = "暂时就用完全版"
$ws.Range("E39").Value = "Y"
$ws.Range("E41").Value = "Y"
$ws.Range("E42").Value = "Y"
$ws.Range("E44").Value = "Y"
$ws.Range("E45").Value = "Y"
$ws.Range("E46").Value = "Y"
$ws.Range("E47").Value = "Y"
$ws.Range("E48").Value = "Y"
$ws.Range("E49").Value = "Y"
$ws.Range("E50").Value = "Y"
$ws.Range("E51").Value = "Y"
$ws.Range("E52").Value = "Y"
$ws.Range("E54").Value = "Y"
$ws.Range("E55").Value = "Y"
$ws.Range("E56").Value = "Y"
$ws.Range("E58").Value = "Y"
$ws.Range("E59").Value = "Y"
$ws.Range("E61").Value = "Y"
$ws.Range("E62").Value = "Y"
$ws.Range("E63").Value = "Y"
$ws.Range("E64").Value = "Y"
$ws.Range("E66").Value = "Y"
$ws.Range("E67").Value = "Y"
$ws.Range("E68").Value = "Y"
$ws.Range("D69").Value = "缺一个？"
$ws.Range("E69").Value = "Y"
$ws.Range("F69").Value = "绿色的密码锁"
$ws.Range("E70").Value = "Y"
$ws.Range("E72").Value = "Y"
$ws.Range("E73").Value = "Y"
$ws.Range("E75").Value = "Y"
$ws.Range("F75").Value = "换成叉烧包OK"
$ws.Range("E76").Value = "Y"
$ws.Range("E77").Value = "Y"

# Column D width
$ws.Columns.Item(4).ColumnWidth = 33.857142857142854

# Update selection to match the final view state
$ws.Range("E78").Select()
